$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 261.57144
$ws.Range("J6").Value = 148.33333
$ws.Range("L6").Value = 444.99999
$ws.Range("N6").Value = -668.99999
$ws.Range("H19").Value = 2302.4
$ws.Range("I19").Value = 1989.25
$ws.Range("J19").Value = 2511.1667
$ws.Range("K19").Value = 1989.25
$ws.Range("L19").Value = 2511.1667
$ws.Range("M19").Value = -1814.25
$ws.Range("N19").Value = -2861.1667
$ws.Range("H28").Value = 768.6667
$ws.Range("I28").Value = 630.2857
$ws.Range("K28").Value = 630.2857
$ws.Range("M28").Value = -145.2857
$ws.Range("H41").Value = 1499.25
$ws.Range("I41").Value = 2109.2
$ws.Range("K41").Value = 2109.2
$ws.Range("M41").Value = -1669.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 23875.334
$ws.Range("J24").Value = 23875.334
$ws.Range("L24").Value = 23875.334
$ws.Range("N24").Value = -24623.334
$ws.Range("H32").Value = 12699.167
$ws.Range("I32").Value = 1084.3
$ws.Range("K32").Value = 1084.3
$ws.Range("M32").Value = -797.3
$ws.Range("H43").Value = 36876.5
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 36876.5
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 36876.5
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -37502.5
$ws.Range("H45").Value = 3255.8696
$ws.Range("I45").Value = 2800.7144
$ws.Range("J45").Value = 3963.889
$ws.Range("K45").Value = 2800.7144
$ws.Range("L45").Value = 3963.889
$ws.Range("M45").Value = -2423.7144
$ws.Range("N45").Value = -4717.889
$ws.Range("H61").Value = 8369.166999999999
$ws.Range("I61").Value = 6043.2
$ws.Range("J61").Value = 19999
$ws.Range("K61").Value = 6043.2
$ws.Range("L61").Value = 19999
$ws.Range("M61").Value = -5831.2
$ws.Range("N61").Value = -20423
$ws.Range("H63").Value = 5137.6
$ws.Range("I63").Value = 5137.6
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 5137.6
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -4451.6
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 5137.6
$ws.Range("I66").Value = 5137.6
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 25688
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -22256
$ws.Range("N66").ClearContents()
$ws.Range("H97").Value = 2155.8696
$ws.Range("I97").Value = 707.5
$ws.Range("K97").Value = 707.5
$ws.Range("M97").Value = -211.5
$ws.Range("H100").Value = 23875.334
$ws.Range("J100").Value = 23875.334
$ws.Range("L100").Value = 23875.334
$ws.Range("N100").Value = -26039.334
$ws.Range("H109").Value = 291351
$ws.Range("J109").Value = 291351
$ws.Range("L109").Value = 291351
$ws.Range("N109").Value = -294125
$ws.Range("H117").Value = 59000
$ws.Range("J117").Value = 59000
$ws.Range("L117").Value = 59000
$ws.Range("N117").Value = -68178
$ws.Range("H136").Value = 8369.166999999999
$ws.Range("I136").Value = 6043.2
$ws.Range("J136").Value = 19999
$ws.Range("K136").Value = 18129.6
$ws.Range("L136").Value = 59997
$ws.Range("M136").Value = -15579.6
$ws.Range("N136").Value = -65097

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2891.3333
$ws.Range("I86").Value = 1749.75
$ws.Range("J86").Value = 5174.5
$ws.Range("K86").Value = 1749.75
$ws.Range("L86").Value = 5174.5
$ws.Range("M86").Value = -626.75
$ws.Range("N86").Value = -7420.5
$ws.Range("H89").Value = 2891.3333
$ws.Range("I89").Value = 1749.75
$ws.Range("J89").Value = 5174.5
$ws.Range("K89").Value = 8748.75
$ws.Range("L89").Value = 25872.5
$ws.Range("M89").Value = -3132.75
$ws.Range("N89").Value = -37104.5
$ws.Range("H134").Value = 2500
$ws.Range("I134").Value = 2500
$ws.Range("K134").Value = 7500
$ws.Range("M134").Value = -4965

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 132437.12
$ws.Range("I16").Value = 55999.6
$ws.Range("J16").Value = 259833
$ws.Range("K16").Value = 55999.6
$ws.Range("L16").Value = 259833
$ws.Range("M16").Value = -55712.6
$ws.Range("N16").Value = -260407
$ws.Range("H113").Value = 132437.12
$ws.Range("I113").Value = 55999.6
$ws.Range("J113").Value = 259833
$ws.Range("K113").Value = 55999.6
$ws.Range("L113").Value = 259833
$ws.Range("M113").Value = -53829.6
$ws.Range("N113").Value = -264173

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 338.5
$ws.Range("I18").Value = 338.5
$ws.Range("K18").Value = 1015.5
$ws.Range("M18").Value = -846.5
$ws.Range("H29").Value = 20883.334
$ws.Range("J29").Value = 6275.25
$ws.Range("L29").Value = 18825.75
$ws.Range("N29").Value = -19379.75
$ws.Range("H46").Value = 40069.152
$ws.Range("I46").Value = 1579.95
$ws.Range("J46").Value = 168366.5
$ws.Range("K46").Value = 4739.85
$ws.Range("L46").Value = 505099.5
$ws.Range("M46").Value = -4648.85
$ws.Range("N46").Value = -505281.5
$ws.Range("H50").Value = 570
$ws.Range("I50").Value = 455.9091
$ws.Range("J50").Value = 988.3333
$ws.Range("K50").Value = 1367.7273
$ws.Range("L50").Value = 2964.9999
$ws.Range("M50").Value = -886.7273
$ws.Range("N50").Value = -3926.9999
$ws.Range("H51").Value = 3709.923
$ws.Range("I51").Value = 3185.75
$ws.Range("K51").Value = 9557.25
$ws.Range("M51").Value = -9097.25
$ws.Range("H53").Value = 570
$ws.Range("I53").Value = 455.9091
$ws.Range("J53").Value = 988.3333
$ws.Range("K53").Value = 1367.7273
$ws.Range("L53").Value = 2964.9999
$ws.Range("M53").Value = -886.7273
$ws.Range("N53").Value = -3926.9999
$ws.Range("H58").Value = 7498.75
$ws.Range("J58").Value = 8331.666999999999
$ws.Range("L58").Value = 24995.001
$ws.Range("N58").Value = -25251.001
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3145.762
$ws.Range("I132").Value = 3427.6
$ws.Range("J132").Value = 2441.1667
$ws.Range("K132").Value = 10282.8
$ws.Range("L132").Value = 7323.500100000001
$ws.Range("M132").Value = -7752.799999999999
$ws.Range("N132").Value = -12383.5001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 22903.824
$ws.Range("I7").Value = 18028.154
$ws.Range("K7").Value = 18028.154
$ws.Range("M7").Value = -17916.154
$ws.Range("H16").Value = 1651.9375
$ws.Range("I16").Value = 1102.0714
$ws.Range("J16").Value = 5501
$ws.Range("K16").Value = 1102.0714
$ws.Range("L16").Value = 5501
$ws.Range("M16").Value = -932.0714
$ws.Range("N16").Value = -5841
$ws.Range("H22").Value = 2148.9092
$ws.Range("I22").Value = 2071.0667
$ws.Range("J22").Value = 2213.7778
$ws.Range("K22").Value = 2071.0667
$ws.Range("L22").Value = 2213.7778
$ws.Range("M22").Value = -1776.0667
$ws.Range("N22").Value = -2803.7778
$ws.Range("H27").Value = 2148.9092
$ws.Range("I27").Value = 2071.0667
$ws.Range("J27").Value = 2213.7778
$ws.Range("K27").Value = 2071.0667
$ws.Range("L27").Value = 2213.7778
$ws.Range("M27").Value = -1964.0667
$ws.Range("N27").Value = -2427.7778
$ws.Range("H46").Value = 4879.4
$ws.Range("I46").Value = 4866
$ws.Range("J46").Value = 4899.5
$ws.Range("K46").Value = 4866
$ws.Range("L46").Value = 4899.5
$ws.Range("M46").Value = -4678
$ws.Range("N46").Value = -5275.5
$ws.Range("H82").Value = 2020.8572
$ws.Range("I82").Value = 2020.8572
$ws.Range("K82").Value = 2020.8572
$ws.Range("M82").Value = -1659.8572
$ws.Range("H85").Value = 2020.8572
$ws.Range("I85").Value = 2020.8572
$ws.Range("K85").Value = 2020.8572
$ws.Range("M85").Value = -772.8571999999999
$ws.Range("H126").Value = 22903.824
$ws.Range("I126").Value = 18028.154
$ws.Range("K126").Value = 54084.462
$ws.Range("M126").Value = -51614.462

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 4400
$ws.Range("I6").Value = 5900
$ws.Range("J6").Value = 1400
$ws.Range("K6").Value = 5900
$ws.Range("L6").Value = 1400
$ws.Range("M6").Value = -5785
$ws.Range("N6").Value = -1630
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("H17").Value = 10000
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 10000
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 10000
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -10344
$ws.Range("H81").Value = 6056.625
$ws.Range("I81").Value = 6379.6
$ws.Range("J81").Value = 5518.3335
$ws.Range("K81").Value = 12759.2
$ws.Range("L81").Value = 11036.667
$ws.Range("M81").Value = -11698.2
$ws.Range("N81").Value = -13158.667
$ws.Range("H84").Value = 6056.625
$ws.Range("I84").Value = 6379.6
$ws.Range("J84").Value = 5518.3335
$ws.Range("K84").Value = 63796
$ws.Range("L84").Value = 55183.335
$ws.Range("M84").Value = -58492
$ws.Range("N84").Value = -65791.33499999999
$ws.Range("H96").Value = 1639.8
$ws.Range("I96").Value = 1600
$ws.Range("J96").Value = 1699.5
$ws.Range("K96").Value = 1600
$ws.Range("L96").Value = 1699.5
$ws.Range("M96").Value = -227
$ws.Range("N96").Value = -4445.5
